$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for rows 2-7, columns B, C, E, F (column D unchanged)
$data = @(
    @{ Row = 2; B = 8766.99705461952;  C = 8461.79515861177;  E = 4143.77854065367; F = -19.9344291972736 },
    @{ Row = 3; B = 8446.95309152382;  C = 7495.01618535487;  E = 4011.99833796952; F = 104.292271805183  },
    @{ Row = 4; B = 2998.14757781871;  C = 5095.32457031298;  E = 4190.6520784596;  F = 11.7490270321906  },
    @{ Row = 5; B = 2815.93090957117;  C = 5049.25668759568;  E = 4128.92726094809; F = 7.2576645226574   },
    @{ Row = 6; B = 9081.33564229631;  C = 8401.49146369489;  E = 5123.12245024051; F = 188.358913080641  },
    @{ Row = 7; B = 9049.92958844127;  C = 8786.89211347041;  E = 5104.08460731287; F = 203.624030032637  }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
}
